$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "49.117.05"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.22%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.628.27"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.45%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "111.63"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.39%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "323.01"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.44%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -1.73%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.09%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.542"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -3.31%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.74"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -2.80%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "19.78"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -4.32%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0810"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -1.45%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.125"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.20%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.25"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.85%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.045.31"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.62%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.620.48"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.35%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.858"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.79%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "49.124.49"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.12%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -4.13%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.89"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -3.18%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.69"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -1.15%  "
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.87%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "269.26"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -4.37%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "68.46"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -5.91%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -1.39%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.09"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -2.03%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.10%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.08"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.93%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.43%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "35.03"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -3.20%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -4.05%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "49.47"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.25%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.23%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.13%  "
$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0797"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.06%  "
$ws.Range("B36").Value = "Celestia"
$ws.Range("C36").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.99"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -3.35%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.97"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +4.40%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.74%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.12"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.96%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "127.51"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +3.55%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.110"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -1.83%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "22.08"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -2.84%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -4.34%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.51%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.064.14"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.54%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -3.49%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +6.10%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -3.69%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.91"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -1.24%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.19"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -3.36%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "58.60"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +1.47%  "
